# Update the crypto price/volume table with the latest scraped values.
# (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Cells.Item(2, 4).Value = '69.761.40'
$ws.Cells.Item(2, 5).Value = '  -0.98%  '

# Row 3: Ethereum
$ws.Cells.Item(3, 4).Value = '3.492.20'
$ws.Cells.Item(3, 5).Value = '  -2.19%  '

# Row 4: TetherUSD
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '0.999'
$ws.Cells.Item(4, 5).Value = '  -0.07%  '

# Row 5: BNB
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '607.45'
$ws.Cells.Item(5, 5).Value = '  +0.46%  '

# Row 6: Solana
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '194.47'
$ws.Cells.Item(6, 5).Value = '  +3.33%  '

# Row 7: XRP
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.625'
$ws.Cells.Item(7, 5).Value = '  +0.36%  '

# Row 8: USDC
$ws.Cells.Item(8, 5).Value = '  -0.12%  '

# Row 9: Dogecoin
$ws.Cells.Item(9, 5).Value = '  -1.82%  '

# Row 10: Cardano
$ws.Cells.Item(10, 5).Value = '  +1.24%  '

# Row 11: Avalanche
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '53.47'
$ws.Cells.Item(11, 5).Value = '  -1.69%  '

# Row 12: ShibaInu
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.0000306'
$ws.Cells.Item(12, 5).Value = '  -2.08%  '

# Row 13: Polkadot
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '9.61'
$ws.Cells.Item(13, 5).Value = '  +1.59%  '

# Row 14: WrappedliquidstakedEther2.0
$ws.Cells.Item(14, 4).Value = '4.058.11'
$ws.Cells.Item(14, 5).Value = '  -1.68%  '

# Row 15: BitcoinCash
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '601.25'
$ws.Cells.Item(15, 5).Value = '  +4.18%  '

# Row 16: WrappedBTC
$ws.Cells.Item(16, 4).Value = '69.856.08'
$ws.Cells.Item(16, 5).Value = '  -0.88%  '

# Row 17: Uniswap
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '12.64'
$ws.Cells.Item(17, 5).Value = '  -1.51%  '

# Row 18: Chainlink
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '18.90'
$ws.Cells.Item(18, 5).Value = '  -0.95%  '

# Row 19: WrappedEther
$ws.Cells.Item(19, 4).Value = '3.505.85'
$ws.Cells.Item(19, 5).Value = '  -2.84%  '

# Row 20: TRON
$ws.Cells.Item(20, 5).Value = '  -0.22%  '

# Row 21: Polygon
$ws.Cells.Item(21, 5).Value = '  -1.19%  '

# Row 22: InternetComputer(DFINITY)
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '17.84'
$ws.Cells.Item(22, 5).Value = '  +1.06%  '

# Row 23: Litecoin
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '104.98'
$ws.Cells.Item(23, 5).Value = '  +10.96%  '

# Row 24: PancakeSwap
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '4.64'
$ws.Cells.Item(24, 5).Value = '  -1.47%  '

# Row 25: Toncoin
$ws.Cells.Item(25, 5).Value = '  +3.33%  '

# Row 26: ImmutableX
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '3.07'
$ws.Cells.Item(26, 5).Value = '  +4.46%  '

# Row 27: RenderToken
$ws.Cells.Item(27, 5).Value = '  -0.55%  '

# Row 28: Filecoin
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '9.84'
$ws.Cells.Item(28, 5).Value = '  +4.53%  '

# Row 29: EthereumClassic
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '34.00'
$ws.Cells.Item(29, 5).Value = '  +4.80%  '

# Row 30: dogwifhat
$ws.Cells.Item(30, 2).Value = 'dogwifhat'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '4.42'
$ws.Cells.Item(30, 5).Value = '  +16.75%  '

# Row 31: NEARProtocol
$ws.Cells.Item(31, 2).Value = 'NEARProtocol'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '7.16'
$ws.Cells.Item(31, 5).Value = '  +0.71%  '

# Row 32: Cosmos
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '12.65'
$ws.Cells.Item(32, 5).Value = '  +3.21%  '

# Row 33: Hedera
$ws.Cells.Item(33, 5).Value = '  +0.29%  '

# Row 34: OKB
$ws.Cells.Item(34, 5).Value = '  +0.27%  '

# Row 35: Maker
$ws.Cells.Item(35, 4).Value = '3.692.99'
$ws.Cells.Item(35, 5).Value = '  -4.35%  '

# Row 36: Dai
$ws.Cells.Item(36, 5).Value = '  -0.11%  '

# Row 37: Bittensor
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '519.23'
$ws.Cells.Item(37, 5).Value = '  -0.49%  '

# Row 38: Fetch.AI
$ws.Cells.Item(38, 5).Value = '  -6.00%  '

# Row 39: PEPE
$ws.Cells.Item(39, 4).Value = '0.0₃0786'
$ws.Cells.Item(39, 5).Value = '  -0.65%  '

# Row 40: TheGraph
$ws.Cells.Item(40, 5).Value = '  -4.34%  '

# Row 41: InjectiveProtocol
$ws.Cells.Item(41, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '36.74'
$ws.Cells.Item(41, 5).Value = '  -3.97%  '

# Row 42: Stacks
$ws.Cells.Item(42, 2).Value = 'Stacks'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '3.57'
$ws.Cells.Item(42, 5).Value = '  +0.13%  '

# Row 43: Kaspa
$ws.Cells.Item(43, 5).Value = '  -1.45%  '

# Row 44: VeChain
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.0460'
$ws.Cells.Item(44, 5).Value = '  +0.81%  '

# Row 45: ThetaToken
$ws.Cells.Item(45, 5).Value = '  -4.24%  '

# Row 46: Stellar
$ws.Cells.Item(46, 5).Value = '  +1.15%  '

# Row 47: ApeXProtocol
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '3.31'
$ws.Cells.Item(47, 5).Value = '  -3.85%  '

# Row 48: THORChain
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '8.77'
$ws.Cells.Item(48, 5).Value = '  -4.96%  '

# Row 49: FirstDigitalUSD
$ws.Cells.Item(49, 5).Value = '  +0.43%  '

# Row 50: Monero
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '132.57'
$ws.Cells.Item(50, 5).Value = '  -2.13%  '

# Row 51: Mantle
$ws.Cells.Item(51, 5).Value = '  +10.34%  '
